# ---------------------------------------------------------------------------
# Edit summary (per the author's commit):
#   1. The table on slide 16 gets a different built-in table style applied
#      (tableStyleId {6DEEC154-FA5E-410C-BF9C-1B9866E95B79} ->
#                    {E745DBA5-AB18-43C6-A8D6-3C2CA43B1A95}).
#   2. The deck's theme colors are switched from the custom "Integral"
#      palette over to the default Office palette (dk1/lt1/dk2/lt2/accent1-6/
#      hlink/folHlink), i.e. a new Design/Theme was applied to the
#      presentation.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Re-style the results table on slide 16 -----------------------------
$oldStyleId = "{6DEEC154-FA5E-410C-BF9C-1B9866E95B79}"
$newStyleId = "{E745DBA5-AB18-43C6-A8D6-3C2CA43B1A95}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2) Apply the Office Theme color palette to the presentation theme -----
function BGR($r, $g, $b) { return $b * 65536 + $g * 256 + $r }

$themeSlide = $p.Slides.Item(1)
$colors = $themeSlide.ThemeColorScheme

$colors.Item(1).RGB  = BGR 0x00 0x00 0x00   # dk1
$colors.Item(2).RGB  = BGR 0xFF 0xFF 0xFF   # lt1
$colors.Item(3).RGB  = BGR 0x44 0x54 0x6A   # dk2
$colors.Item(4).RGB  = BGR 0xE7 0xE6 0xE6   # lt2
$colors.Item(5).RGB  = BGR 0x5B 0x9B 0xD5   # accent1
$colors.Item(6).RGB  = BGR 0xED 0x7D 0x31   # accent2
$colors.Item(7).RGB  = BGR 0xA5 0xA5 0xA5   # accent3
$colors.Item(8).RGB  = BGR 0xFF 0xC0 0x00   # accent4
$colors.Item(9).RGB  = BGR 0x44 0x72 0xC4   # accent5
$colors.Item(10).RGB = BGR 0x70 0xAD 0x47   # accent6
$colors.Item(11).RGB = BGR 0x05 0x63 0xC1   # hlink
$colors.Item(12).RGB = BGR 0x95 0x4F 0x72   # folHlink
